$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.565802097320557
$ws.Range("B1").Value = 2.702942371368408
$ws.Range("C1").Value = 3.072072267532349
$ws.Range("D1").Value = 2.874539613723755
$ws.Range("E1").Value = 3.051093816757202
